# Updated for EA 23.252 Stable: add the "Night Vision At Night" toggle/tooltip rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: toggle07
$ws.Range("A16").Value = "toggle07"
$ws.Range("B16").Value = "Enable Night Vision At Night"
$ws.Range("C16").Value = "夜間に暗視を有効にする"
$ws.Range("D16").Value = "启用夜间夜视"

# Row 17: tooltip07
$ws.Range("A17").Value = "tooltip07"
$ws.Range("B17").Value = "Grants Night Vision (Cat's Eye) automatically while it is night."
$ws.Range("C17").Value = "夜間になると自動的に暗視（猫の目）を付与します。"
$ws.Range("D17").Value = "在夜间自动赋予夜视（猫之眼）效果。"

# Match the wrapped "Noto Sans SC" formatting used by the other jp/cn tooltip & label cells.
$wrapRange = $ws.Range("C16:D17")
$wrapRange.Font.Name = "Noto Sans SC"
$wrapRange.WrapText = $true

# Restore the view to the top-left corner and leave the active cell where the author left it.
$ws.Range("D23").Select()
